$wb = $excel.ActiveWorkbook

# --- Swap the names of the "bubble"/"out" playoff sheets (trad + adv) ---
# Sheet positions/ids/rIds stay put; only the display names trade places.
$wsBubbleTrad = $wb.Worksheets.Item("bubbleplayofftrad")
$wsOutTrad    = $wb.Worksheets.Item("outplayofftrad")
$wsBubbleTrad.Name = "__swap_trad__"
$wsOutTrad.Name    = "bubbleplayofftrad"
$wsBubbleTrad.Name = "outplayofftrad"

$wsBubbleAdv = $wb.Worksheets.Item("bubbleplayoffadv")
$wsOutAdv    = $wb.Worksheets.Item("outplayoffadv")
$wsBubbleAdv.Name = "__swap_adv__"
$wsOutAdv.Name    = "bubbleplayoffadv"
$wsBubbleAdv.Name = "outplayoffadv"

# --- Update per-sheet selections ---
# The sheet that used to be "outplayofftrad" (now "bubbleplayofftrad") keeps
# its own cell selection (I21) but is no longer the active/selected tab.
[void]$wsOutTrad.Range("I21").Select()

# The sheet that used to be "bubbleplayoffadv" (now "outplayoffadv") becomes
# the active tab, with a new selection of M28.
[void]$wsBubbleAdv.Activate()
[void]$wsBubbleAdv.Range("M28").Select()
